$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46 (gftfpp / nonprofit_ppp) for the new
# "gftfpr" / "medicare_reimbursement_increase" (Medicare reimbursement increase) code.
$ws.Rows.Item(46).Insert()
$ws.Range("A46").Value = "gftfpr"
$ws.Range("B46").Value = "medicare_reimbursement_increase"

# Insert a new row before the "gfegv" / "provider_relief_fund" row (now row 57)
# for the new "gfsubv" / "provider_relief_fund" code.
$ws.Rows.Item(54).Insert()
$ws.Range("A54").Value = "gfsubv"
$ws.Range("B54").Value = "provider_relief_fund"

# Update the existing "gfegv" row's reference to the more specific name.
$ws.Range("B58").Value = "provider_relief_fund_grants"
